$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-16 Wednesday" "2024-10-17 Thursday"

Replace-Text "217×7=" "496×4="
Replace-Text "682×8=" "736×4="
Replace-Text "239×8=" "573×3="
Replace-Text "654×5=" "163×8="
Replace-Text "460×7=" "748×3="

Replace-Text "279×4=" "434×3="
Replace-Text "210×6=" "859×6="
Replace-Text "384×3=" "644×8="
Replace-Text "367×6=" "677×2="
Replace-Text "366×4=" "784×6="

Replace-Text "157×7=" "408×3="
Replace-Text "674×6=" "342×2="
Replace-Text "386×9=" "516×2="
Replace-Text "351×5=" "233×9="
Replace-Text "761×3=" "833×4="

Replace-Text "740×6=" "521×5="
Replace-Text "159×6=" "419×4="
Replace-Text "842×5=" "237×8="
Replace-Text "818×3=" "561×9="
Replace-Text "968×5=" "185×6="

Replace-Text "879×2=" "506×5="
Replace-Text "466×2=" "632×6="
Replace-Text "214×6=" "145×3="
Replace-Text "180×3=" "260×3="
Replace-Text "205×7=" "144×8="
